$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 4355.8887
$ws.Range("I113").Value = 4034
$ws.Range("J113").Value = 4999.6665
$ws.Range("K113").Value = 4034
$ws.Range("L113").Value = 4999.6665
$ws.Range("M113").Value = -780
$ws.Range("N113").Value = -11507.6665

$ws.Range("H137").Value = 23763.334
$ws.Range("I137").Value = 23763.334
$ws.Range("K137").Value = 71290.00199999999
$ws.Range("M137").Value = -68740.00199999999

$ws.Range("H138").Value = 2314.923
$ws.Range("J138").Value = 3070.8333
$ws.Range("L138").Value = 9212.499899999999
$ws.Range("N138").Value = -19492.4999

$ws.Range("H141").Value = 6261
$ws.Range("I141").Value = 3138.111
$ws.Range("K141").Value = 9414.332999999999
$ws.Range("M141").Value = -4234.332999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2106
$ws.Range("I2").Value = 942
$ws.Range("J2").Value = 3270
$ws.Range("K2").Value = 942
$ws.Range("L2").Value = 3270
$ws.Range("M2").Value = -829
$ws.Range("N2").Value = -3496

$ws.Range("H32").Value = 2908.8723
$ws.Range("I32").Value = 2958.761
$ws.Range("K32").Value = 2958.761
$ws.Range("M32").Value = -2671.761

$ws.Range("H45").Value = 768.25
$ws.Range("I45").Value = 768.25
$ws.Range("K45").Value = 768.25
$ws.Range("M45").Value = -391.25

$ws.Range("H74").Value = 2408.0476
$ws.Range("I74").Value = 2408.0476
$ws.Range("K74").Value = 2408.0476
$ws.Range("M74").Value = -1534.0476

$ws.Range("H77").Value = 2408.0476
$ws.Range("I77").Value = 2408.0476
$ws.Range("K77").Value = 12040.238
$ws.Range("M77").Value = -7672.237999999999

$ws.Range("H110").Value = 994.7
$ws.Range("I110").Value = 994.7
$ws.Range("K110").Value = 994.7
$ws.Range("M110").Value = 1050.3

$ws.Range("H116").Value = 2106
$ws.Range("I116").Value = 942
$ws.Range("J116").Value = 3270
$ws.Range("K116").Value = 942
$ws.Range("L116").Value = 3270
$ws.Range("M116").Value = 1352
$ws.Range("N116").Value = -7858

$ws.Range("H122").Value = 2254.4285
$ws.Range("I122").Value = 1799.1765
$ws.Range("J122").Value = 2958
$ws.Range("K122").Value = 5397.529500000001
$ws.Range("L122").Value = 8874
$ws.Range("M122").Value = -2947.529500000001
$ws.Range("N122").Value = -13774

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2106
$ws.Range("I3").Value = 942
$ws.Range("J3").Value = 3270
$ws.Range("K3").Value = 942
$ws.Range("L3").Value = 3270
$ws.Range("M3").Value = -828
$ws.Range("N3").Value = -3498

$ws.Range("H20").Value = 1586.5
$ws.Range("I20").Value = 1586.5
$ws.Range("K20").Value = 1586.5
$ws.Range("M20").Value = -1339.5

$ws.Range("H22").Value = 857.76
$ws.Range("I22").Value = 668.93335
$ws.Range("K22").Value = 668.93335
$ws.Range("M22").Value = -495.93335

$ws.Range("H26").Value = 45498.168
$ws.Range("I26").Value = 38799.8
$ws.Range("K26").Value = 38799.8
$ws.Range("M26").Value = -38507.8

$ws.Range("H94").Value = 2818
$ws.Range("J94").Value = 2739.8
$ws.Range("L94").Value = 2739.8
$ws.Range("N94").Value = -3641.8

$ws.Range("H105").Value = 1891.6666
$ws.Range("I105").Value = 1579.1666
$ws.Range("K105").Value = 1579.1666
$ws.Range("M105").Value = 167.8334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 2265.4
$ws.Range("I2").Value = 2668.625
$ws.Range("J2").Value = 652.5
$ws.Range("K2").Value = 2668.625
$ws.Range("L2").Value = 652.5
$ws.Range("M2").Value = -2555.625
$ws.Range("N2").Value = -878.5

$ws.Range("H5").Value = 361.3
$ws.Range("I5").Value = 290.33334
$ws.Range("J5").Value = 1000
$ws.Range("K5").Value = 290.33334
$ws.Range("L5").Value = 1000
$ws.Range("M5").Value = -178.33334
$ws.Range("N5").Value = -1224

$ws.Range("H16").Value = 3105.3333
$ws.Range("I16").Value = 2908.25
$ws.Range("K16").Value = 2908.25
$ws.Range("M16").Value = -2621.25

$ws.Range("H17").Value = 9545.454
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 9545.454
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 9545.454
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -9893.454

$ws.Range("H19").Value = 2412.5
$ws.Range("I19").Value = 2263.75
$ws.Range("J19").Value = 2710
$ws.Range("K19").Value = 2263.75
$ws.Range("L19").Value = 2710
$ws.Range("M19").Value = -2093.75
$ws.Range("N19").Value = -3050

$ws.Range("H24").Value = 2412.5
$ws.Range("I24").Value = 2263.75
$ws.Range("J24").Value = 2710
$ws.Range("K24").Value = 2263.75
$ws.Range("L24").Value = 2710
$ws.Range("M24").Value = -2093.75
$ws.Range("N24").Value = -3050

$ws.Range("H31").Value = 4485.4814
$ws.Range("I31").Value = 1785.5714
$ws.Range("J31").Value = 13935.167
$ws.Range("K31").Value = 1785.5714
$ws.Range("L31").Value = 13935.167
$ws.Range("M31").Value = -1490.5714
$ws.Range("N31").Value = -14525.167

$ws.Range("H34").Value = 4485.4814
$ws.Range("I34").Value = 1785.5714
$ws.Range("J34").Value = 13935.167
$ws.Range("K34").Value = 1785.5714
$ws.Range("L34").Value = 13935.167
$ws.Range("M34").Value = -1583.5714
$ws.Range("N34").Value = -14339.167

$ws.Range("H99").Value = 3085.5715
$ws.Range("I99").Value = 2780
$ws.Range("J99").Value = 3849.5
$ws.Range("K99").Value = 2780
$ws.Range("L99").Value = 3849.5
$ws.Range("M99").Value = -1282
$ws.Range("N99").Value = -6845.5

$ws.Range("H105").Value = 1244.4117
$ws.Range("I105").Value = 953.1429000000001
$ws.Range("J105").Value = 2603.6667
$ws.Range("K105").Value = 953.1429000000001
$ws.Range("L105").Value = 2603.6667
$ws.Range("M105").Value = 793.8570999999999
$ws.Range("N105").Value = -6097.6667

$ws.Range("H113").Value = 3105.3333
$ws.Range("I113").Value = 2908.25
$ws.Range("K113").Value = 2908.25
$ws.Range("M113").Value = -738.25

$ws.Range("H122").Value = 1869.6154
$ws.Range("I122").Value = 1692.0834
$ws.Range("K122").Value = 5076.2502
$ws.Range("M122").Value = -2626.2502

$ws.Range("H126").Value = 3085.5715
$ws.Range("I126").Value = 2780
$ws.Range("J126").Value = 3849.5
$ws.Range("K126").Value = 8340
$ws.Range("L126").Value = 11548.5
$ws.Range("M126").Value = -5870
$ws.Range("N126").Value = -16488.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 2795.6
$ws.Range("I3").Value = 2795.6
$ws.Range("K3").Value = 8386.799999999999
$ws.Range("M3").Value = -8274.799999999999

$ws.Range("H4").Value = 97990340
$ws.Range("I4").Value = 122987570
$ws.Range("J4").Value = 10500000
$ws.Range("K4").Value = 368962710
$ws.Range("L4").Value = 31500000
$ws.Range("M4").Value = -368962598
$ws.Range("N4").Value = -31500224

$ws.Range("H32").Value = 10100
$ws.Range("J32").Value = 10100
$ws.Range("L32").Value = 30300
$ws.Range("N32").Value = -30866

$ws.Range("H33").Value = 136.7
$ws.Range("I33").Value = 110.666664
$ws.Range("K33").Value = 663.999984
$ws.Range("M33").Value = -380.999984

$ws.Range("H55").Value = 2091014
$ws.Range("I55").Value = 1000432.6
$ws.Range("J55").Value = 2870000.8
$ws.Range("K55").Value = 3001297.8
$ws.Range("L55").Value = 8610002.399999999
$ws.Range("M55").Value = -3001120.8
$ws.Range("N55").Value = -8610356.399999999

$ws.Range("H107").Value = 2186
$ws.Range("J107").Value = 944.7778
$ws.Range("L107").Value = 2834.3334
$ws.Range("N107").Value = -6674.3334

$ws.Range("H121").Value = 655.2857
$ws.Range("I121").Value = 520
$ws.Range("K121").Value = 1560
$ws.Range("M121").Value = -250

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1559.1538
$ws.Range("J97").Value = 1554.1666
$ws.Range("L97").Value = 1554.1666
$ws.Range("N97").Value = -2546.1666

$ws.Range("H132").Value = 3811.8667
$ws.Range("I132").Value = 3597.5833
$ws.Range("J132").Value = 4669
$ws.Range("K132").Value = 10792.7499
$ws.Range("L132").Value = 14007
$ws.Range("M132").Value = -8262.749899999999
$ws.Range("N132").Value = -19067

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 6055.391
$ws.Range("I93").Value = 4795.364
$ws.Range("K93").Value = 4795.364
$ws.Range("M93").Value = -3547.364

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 10000
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 10000
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 10000
$ws.Range("M13").ClearContents()
$ws.Range("N13").Value = -10280

$ws.Range("H62").Value = 11750.167
$ws.Range("J62").Value = 11818.182
$ws.Range("L62").Value = 11818.182
$ws.Range("N62").Value = -13066.182

$ws.Range("H65").Value = 11750.167
$ws.Range("J65").Value = 11818.182
$ws.Range("L65").Value = 59090.91
$ws.Range("N65").Value = -65330.91

$ws.Range("H81").Value = 1366.9231
$ws.Range("I81").Value = 1179.1111
$ws.Range("J81").Value = 1789.5
$ws.Range("K81").Value = 2358.2222
$ws.Range("L81").Value = 3579
$ws.Range("M81").Value = -1297.2222
$ws.Range("N81").Value = -5701

$ws.Range("H84").Value = 1366.9231
$ws.Range("I84").Value = 1179.1111
$ws.Range("J84").Value = 1789.5
$ws.Range("K84").Value = 11791.111
$ws.Range("L84").Value = 17895
$ws.Range("M84").Value = -6487.111000000001
$ws.Range("N84").Value = -28503

$ws.Range("H100").Value = 884.2857
$ws.Range("I100").Value = 908
$ws.Range("K100").Value = 1816
$ws.Range("M100").Value = -1275

$ws.Range("H126").Value = 1633.8
$ws.Range("I126").Value = 1633.8
$ws.Range("K126").Value = 4901.4
$ws.Range("M126").Value = -2431.4

$ws.Range("H133").Value = 77999.664
$ws.Range("J133").Value = 77999.664
$ws.Range("L133").Value = 77999.664
$ws.Range("N133").Value = -88119.664
